$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 12.3292
$ws.Range("A8").Value = -21.1186
$ws.Range("A10").Value = -20.49849999999997
$ws.Range("A12").Value = -22.36820000000003
$ws.Range("D13").Value = -7.836800000000002
$ws.Range("A18").Value = -22.25920000000003
$ws.Range("E20").Value = 12.16589999999999
